$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; F=1; G=2.021590333333334; H=6.064771; I=0.01116262347650641; J=0.01116262347650641; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4475167411500002; P=0.4475167411500002; Q=8.567177052076779; R=77.10459346869101; S=0.004995460880890635; T=0.004995460880890635 }
    3 = @{ E=3; F=1; G=2.021590333333334; H=6.064771; I=0.01116262347650641; J=0.01116262347650641; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5524832588499998; P=0.5524832588499998; Q=10.57663649568322; R=95.189728461149; S=0.006167162595615778; T=0.006167162595615778 }
    4 = @{ E=3; F=1; G=25.140634; H=75.421902; I=0.1388191398995883; J=0.1388191398995883; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4475167411500002; P=0.4475167411500002; Q=106.5419927707713; R=958.8779349369421; S=0.06212388909710972; T=0.0621238890971097 }
    5 = @{ E=3; F=1; G=25.140634; H=75.421902; I=0.1388191398995883; J=0.1388191398995883; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5524832588499998; P=0.5524832588499998; Q=131.5317662063487; R=1183.785895857138; S=0.07669525080247858; T=0.07669525080247856 }
    6 = @{ E=3; F=1; G=153.9412893333333; H=461.823868; I=0.8500182366239053; J=0.8500182366239052; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4475167411500002; P=0.4475167411500002; Q=652.3786049021365; R=5871.407444119228; S=0.3803973911719998; T=0.3803973911719998 }
    7 = @{ E=3; F=1; G=153.9412893333333; H=461.823868; I=0.8500182366239053; J=0.8500182366239052; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5524832588499998; P=0.5524832588499998; Q=805.3961438719437; R=7248.565294847492; S=0.4696208454519055; T=0.4696208454519054 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
